# support nested json, support datetime as iso8601
#
# Adds a new "time" row (Sheet1 row 15) above the existing "command_test"
# row, shifting command_test / v.e.r.y...deep.path / test / otherprefix
# rows down by one. Updates the json.* defined names to match, adds two
# new defined names (json.datetime, json.time), and moves the sheet
# selection to A15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Insert a new row above the current row 15 ("command_test"), pushing
#    command_test/deep-path/test/otherprefix rows down by one.
$ws.Rows.Item(15).Insert()

# 2. Populate the new row 15 with the "time" label + a time-of-day value.
$ws.Range("A15").Value = "time"
$ws.Range("J15").Value = 0.52425925925925931
$ws.Range("J15").NumberFormat = "[`$-F400]h:mm:ss\ AM/PM"

# 3. Fix up the defined names that pointed at rows which shifted down.
$wb.Names.Item("json.command_test").RefersTo = "=Sheet1!`$J`$16"
$wb.Names.Item("json.v.e.r.y.v.e.r.y.d.e.e.p.p.a.t.h").RefersTo = "=Sheet1!`$J`$17"
$wb.Names.Item("otherprefix.test").RefersTo = "=Sheet1!`$J`$20"

# 4. Register the two new named ranges for the nested/datetime support.
$wb.Names.Add("json.datetime", "=Sheet1!`$J`$14")
$wb.Names.Add("json.time", "=Sheet1!`$J`$15")

# 5. Move the active selection onto the new row (matches author's edit).
$ws.Range("A15").Select()
